$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E6").Value = 3.188
$ws.Range("F6").Value = 3.737
$ws.Range("G6").Value = 2.878
$ws.Range("H6").Value = 3.217
$ws.Range("I6").Value = 2.672
$ws.Range("J6").Value = 2.011

$ws.Range("E39").Value = 6.022
$ws.Range("F39").Value = 6.817
$ws.Range("G39").Value = 5.793
$ws.Range("H39").Value = 5.721
$ws.Range("I39").Value = 4.354
$ws.Range("J39").Value = 4.089
